# Add five new "Day 26" .. "Day 29 employee_uni" worksheets at the end of
# the workbook, populate them with data, and update the selection left
# behind on the previously-last sheet ("Day 25").

$wb = $excel.ActiveWorkbook

function Add-SheetAfterLast {
    param([string]$Name)
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $Name
    return $newSheet
}

# --- "Day 25" used to be the last / active sheet. Excel moves its own
# selection when the user clicks elsewhere before adding more sheets. ---
$day25 = $wb.Worksheets.Item("Day 25")
$day25.Activate()
$day25.Range("I5").Select()

# --- Day 26 : sell_date / product -----------------------------------
$ws = Add-SheetAfterLast "Day 26"

$ws.Range("A1").Value = "sell_date"
$ws.Range("B1").Value = "product"

# Serial day numbers (Excel 1900 date system) for 5/30/2020, 6/1/2020, 6/2/2020
$day26Dates = @(43981, 43983, 43984, 43981, 43983, 43984, 43981)
$day26Products = @("Headphone", "pencil", "Mask", "Basketball", "Bible", "Mask", "T-Shirt")

for ($i = 0; $i -lt $day26Dates.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $day26Dates[$i]
    $ws.Range("A$r").NumberFormat = "m/d/yy"
    $ws.Range("B$r").Value = $day26Products[$i]
}

$ws.Range("G9").Select()

# --- Day 27 : date_id / make_name / lead_id / partner_id -------------
$ws = Add-SheetAfterLast "Day 27"

$ws.Range("A1").Value = "date_id"
$ws.Range("B1").Value = "make_name"
$ws.Range("C1").Value = "lead_id"
$ws.Range("D1").Value = "partner_id"

# Serial day numbers for 12/8/2020 and 12/7/2020
$d1 = 44173
$d2 = 44172

$day27Rows = @(
    @($d1, "toyota", 0, 1),
    @($d1, "toyota", 1, 0),
    @($d1, "toyota", 1, 2),
    @($d2, "toyota", 0, 2),
    @($d2, "toyota", 0, 1),
    @($d1, "honda", 1, 2),
    @($d1, "honda", 2, 1),
    @($d2, "honda", 0, 1),
    @($d2, "honda", 1, 2),
    @($d2, "honda", 2, 1)
)

for ($i = 0; $i -lt $day27Rows.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $day27Rows[$i][0]
    $ws.Range("A$r").NumberFormat = "m/d/yy"
    $ws.Range("B$r").Value = $day27Rows[$i][1]
    $ws.Range("C$r").Value = $day27Rows[$i][2]
    $ws.Range("D$r").Value = $day27Rows[$i][3]
}

$ws.Range("D2").Select()

# --- Day 28 : actor_id / director_id / timestamp ----------------------
$ws = Add-SheetAfterLast "Day 28"

$ws.Range("A1").Value = "actor_id"
$ws.Range("B1").Value = "director_id"
$ws.Range("C1").Value = "timestamp"

$day28Rows = @(
    @(1, 1, 0),
    @(1, 1, 1),
    @(1, 1, 2),
    @(1, 2, 3),
    @(1, 2, 4),
    @(2, 1, 5),
    @(2, 1, 6)
)

for ($i = 0; $i -lt $day28Rows.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $day28Rows[$i][0]
    $ws.Range("B$r").Value = $day28Rows[$i][1]
    $ws.Range("C$r").Value = $day28Rows[$i][2]
}

$ws.Range("E13").Select()

# --- Day 29 employees : id / name -------------------------------------
$ws = Add-SheetAfterLast "Day 29 employees"

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"

$day29EmpRows = @(
    @(1, "Alice"),
    @(7, "Bob"),
    @(11, "Meir"),
    @(90, "Winston"),
    @(3, "Jonathan")
)

for ($i = 0; $i -lt $day29EmpRows.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $day29EmpRows[$i][0]
    $ws.Range("B$r").Value = $day29EmpRows[$i][1]
}

$ws.Range("A1:B6").Select()

# --- Day 29 employee_uni : id / unique_id -----------------------------
$ws = Add-SheetAfterLast "Day 29 employee_uni"

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "unique_id"

$day29UniRows = @(
    @(3, 1),
    @(11, 2),
    @(90, 3)
)

for ($i = 0; $i -lt $day29UniRows.Count; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $day29UniRows[$i][0]
    $ws.Range("B$r").Value = $day29UniRows[$i][1]
}

$ws.Range("L13").Select()
